$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename column headers to short machine-readable codes ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalization fixes of connector words (de/el/los -> De/El/Los) ---
$ws.Range("A9").Value  = "Ciudad De México"
$ws.Range("A13").Value = "Coahuila De Zaragoza"
$ws.Range("A15").Value = "Estado De México"
$ws.Range("B15").Value = "Naucalpan De Juárez"
$ws.Range("B19").Value = "Pachuca De Soto"
$ws.Range("B24").Value = "San Miguel El Alto"
$ws.Range("B26").Value = "Valle De Guadalupe"
$ws.Range("A29").Value = "Michoacán De Ocampo"
$ws.Range("B42").Value = "Mexquitic De Carmona"
$ws.Range("A48").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B48").Value = "Amatlán De Los Reyes"
$ws.Range("B49").Value = "Ignacio De La Llave"

# --- Fix TOTAL -> Total on the grand total row ---
$ws.Range("A56").Value = "Total"

# --- Remove trailing metadata/footer rows (58-62) ---
$ws.Rows("58:62").Delete()
